# A new weekly price observation for "Jengibre" (Terminal La Palmera de La
# Serena) needs to be inserted at the top of the data table (row 62, right
# after the existing rows that stay untouched), pushing all the following
# rows down by one. The new row carries the latest reading; everything below
# it keeps the same relative order it already had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 62..153 down to 63..154, leaving a blank row 62 to fill in.
$ws.Rows("62:62").Insert()

# Fill the newly inserted row with the new observation. The descriptive
# columns (market, region, product, quality, unit, origin, classification)
# are identical for every row in this sheet.
$ws.Cells.Item(62, 1).Value  = 8
$ws.Cells.Item(62, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(62, 3).Value  = "Coquimbo"
$ws.Cells.Item(62, 4).Value2 = 45125
$ws.Cells.Item(62, 5).Value  = 4
$ws.Cells.Item(62, 6).Value  = 100114007
$ws.Cells.Item(62, 7).Value  = "Jengibre"
$ws.Cells.Item(62, 8).Value  = "Sin especificar"
$ws.Cells.Item(62, 9).Value  = "Primera"
$ws.Cells.Item(62, 10).Value = 440
$ws.Cells.Item(62, 11).Value = 17000
$ws.Cells.Item(62, 12).Value = 18000
$ws.Cells.Item(62, 13).Value = 17500
$ws.Cells.Item(62, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(62, 15).Value = "Perú"
$ws.Cells.Item(62, 16).Value = 1346
$ws.Cells.Item(62, 17).Value = 13
$ws.Cells.Item(62, 18).Value = "Hortaliza"
